# feat: add 2022-Q3 data
#
# - "总计" (summary) sheet gains a new 2022-Q3 row (pushing the existing
#   2022-Q1 row down to row 3).
# - A new "2022-Q3" worksheet is inserted between "总计" and "2022-Q1",
#   holding the Q3 fund holding breakdown.
# - The original "2022-Q1" worksheet/data is preserved unchanged, just
#   moved one position later in the tab order.

$xlPasteFormats = -4122

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "总计"
$ws2 = $wb.Worksheets.Item(2)   # currently "2022-Q1"

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift the existing 2022-Q1 total down to row 3, and
#    write the new 2022-Q3 total into row 2.
# ---------------------------------------------------------------------
$ws1.Cells.Item(3, 1).Value = 1
$ws1.Cells.Item(3, 2).Value = "2022-Q1"
$ws1.Cells.Item(3, 3).Value = 1
$ws1.Cells.Item(3, 4).Value = 0.62

# Row 3's first column should look like row 2's (bold/bordered) style.
$ws1.Cells.Item(2, 1).Copy() | Out-Null
$ws1.Cells.Item(3, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws1.Cells.Item(2, 2).Value = "2022-Q3"
$ws1.Cells.Item(2, 4).Value = 0.75

# ---------------------------------------------------------------------
# 2) Add a brand-new worksheet right after the existing "2022-Q1" sheet;
#    it will receive the original, untouched 2022-Q1 fund breakdown.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)

# Carry over the formatting currently on the "2022-Q1" sheet (header row
# + A2) before that sheet gets repurposed below.
$ws2.Range("B1:H1").Copy() | Out-Null
$ws3.Range("B1:H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Cells.Item(2, 1).Copy() | Out-Null
$ws3.Cells.Item(2, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(1, 2).Value = "基金代码"
$ws3.Cells.Item(1, 3).Value = "基金名称"
$ws3.Cells.Item(1, 4).Value = "基金规模"
$ws3.Cells.Item(1, 5).Value = "股票总仓位"
$ws3.Cells.Item(1, 6).Value = "仓位占比"
$ws3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws3.Cells.Item(1, 8).Value = "仓位排名"

$ws3.Cells.Item(2, 1).Value = 0
$ws3.Cells.Item(2, 2).Value = "968029"
$ws3.Cells.Item(2, 3).Value = "恒生指数基金M类人民币（对冲）份额"
$ws3.Cells.Item(2, 4).Value = "25.09"
$ws3.Cells.Item(2, 5).Value = "97.94"
$ws3.Cells.Item(2, 6).Value = "2.47"
$ws3.Cells.Item(2, 7).Value = "0.6197"
$ws3.Cells.Item(2, 8).Value = 10

# ---------------------------------------------------------------------
# 3) Repurpose the original "2022-Q1" sheet object to become the new
#    "2022-Q3" sheet (this keeps the sheet ordering/ids matching a
#    freshly-inserted middle tab: 总计, 2022-Q3, 2022-Q1).
# ---------------------------------------------------------------------
$ws2.Name = "2022-Q3"

# Its header/A2 formatting should match "总计"'s header style.
$ws1.Range("B1:D1").Copy() | Out-Null
$ws2.Range("B1:D1").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("E1:H1").PasteSpecial($xlPasteFormats) | Out-Null
$ws1.Cells.Item(2, 1).Copy() | Out-Null
$ws2.Cells.Item(2, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws2.Cells.Item(2, 4).Value = "27.03"
$ws2.Cells.Item(2, 5).Value = "99.07"
$ws2.Cells.Item(2, 6).Value = "2.77"
$ws2.Cells.Item(2, 7).Value = "0.7487"
$ws2.Cells.Item(2, 8).Value = 9

# ---------------------------------------------------------------------
# 4) Give the newly-added sheet its final name.
# ---------------------------------------------------------------------
$ws3.Name = "2022-Q1"
